# Auto-generated Excel COM-interop script to apply Golem_Profits data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Cells.Item(39, 8).Value = 10
$ws.Cells.Item(39, 9).Value = 10
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 30
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = 266
$ws.Cells.Item(39, 14).ClearContents()
# Row 51
$ws.Cells.Item(51, 8).Value = 142949.75
$ws.Cells.Item(51, 10).Value = 187266.33
$ws.Cells.Item(51, 12).Value = 187266.33
$ws.Cells.Item(51, 14).Value = -188234.33
# Row 64
$ws.Cells.Item(64, 8).Value = 1499.8334
$ws.Cells.Item(64, 9).Value = 1499
$ws.Cells.Item(64, 11).Value = 1499
$ws.Cells.Item(64, 13).Value = -1251
# Row 67
$ws.Cells.Item(67, 8).Value = 1499.8334
$ws.Cells.Item(67, 9).Value = 1499
$ws.Cells.Item(67, 11).Value = 1499
$ws.Cells.Item(67, 13).Value = -641
# Row 97
$ws.Cells.Item(97, 8).Value = 966.6667
$ws.Cells.Item(97, 10).Value = 966.6667
$ws.Cells.Item(97, 12).Value = 2900.0001
$ws.Cells.Item(97, 14).Value = -3892.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Cells.Item(8, 8).Value = 625.5
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 625.5
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 625.5
$ws.Cells.Item(8, 13).ClearContents()
$ws.Cells.Item(8, 14).Value = -913.5

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Cells.Item(2, 8).Value = 2556.2856
$ws.Cells.Item(2, 10).Value = 3298.5
$ws.Cells.Item(2, 12).Value = 3298.5
$ws.Cells.Item(2, 14).Value = -3524.5
# Row 4
$ws.Cells.Item(4, 8).Value = 50001750
$ws.Cells.Item(4, 10).Value = 50001750
$ws.Cells.Item(4, 12).Value = 50001750
$ws.Cells.Item(4, 14).Value = -50001974
# Row 5
$ws.Cells.Item(5, 8).Value = 1931.2858
$ws.Cells.Item(5, 9).Value = 2219.8333
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(5, 11).Value = 2219.8333
$ws.Cells.Item(5, 12).Value = 200
$ws.Cells.Item(5, 13).Value = -2107.8333
$ws.Cells.Item(5, 14).Value = -424
# Row 8
$ws.Cells.Item(8, 8).Value = 4374.75
$ws.Cells.Item(8, 9).Value = 1749.5
$ws.Cells.Item(8, 10).Value = 7000
$ws.Cells.Item(8, 11).Value = 1749.5
$ws.Cells.Item(8, 12).Value = 7000
$ws.Cells.Item(8, 13).Value = -1609.5
$ws.Cells.Item(8, 14).Value = -7280
# Row 10
$ws.Cells.Item(10, 8).Value = 818.6667
$ws.Cells.Item(10, 9).Value = 818.6667
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 818.6667
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -679.6667
$ws.Cells.Item(10, 14).ClearContents()
# Row 11
$ws.Cells.Item(11, 8).Value = 39006
$ws.Cells.Item(11, 10).Value = 39006
$ws.Cells.Item(11, 12).Value = 39006
$ws.Cells.Item(11, 14).Value = -39286
# Row 12
$ws.Cells.Item(12, 8).Value = 672
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(12, 10).Value = 2006
$ws.Cells.Item(12, 11).Value = 5
$ws.Cells.Item(12, 12).Value = 2006
$ws.Cells.Item(12, 13).Value = 165
$ws.Cells.Item(12, 14).Value = -2346
# Row 13
$ws.Cells.Item(13, 8).Value = 309
$ws.Cells.Item(13, 10).Value = 309
$ws.Cells.Item(13, 12).Value = 309
$ws.Cells.Item(13, 14).Value = -587
# Row 14
$ws.Cells.Item(14, 8).Value = 100
$ws.Cells.Item(14, 9).Value = 100
$ws.Cells.Item(14, 11).Value = 100
$ws.Cells.Item(14, 13).Value = 70
# Row 15
$ws.Cells.Item(15, 8).Value = 12926
$ws.Cells.Item(15, 10).Value = 12926
$ws.Cells.Item(15, 12).Value = 12926
$ws.Cells.Item(15, 14).Value = -13266
# Row 19
$ws.Cells.Item(19, 8).Value = 1561.9166
$ws.Cells.Item(19, 9).Value = 158.54546
$ws.Cells.Item(19, 10).Value = 16999
$ws.Cells.Item(19, 11).Value = 158.54546
$ws.Cells.Item(19, 12).Value = 16999
$ws.Cells.Item(19, 13).Value = 11.45454000000001
$ws.Cells.Item(19, 14).Value = -17339
# Row 24
$ws.Cells.Item(24, 8).Value = 1561.9166
$ws.Cells.Item(24, 9).Value = 158.54546
$ws.Cells.Item(24, 10).Value = 16999
$ws.Cells.Item(24, 11).Value = 158.54546
$ws.Cells.Item(24, 12).Value = 16999
$ws.Cells.Item(24, 13).Value = 11.45454000000001
$ws.Cells.Item(24, 14).Value = -17339
# Row 25
$ws.Cells.Item(25, 8).Value = 1000
$ws.Cells.Item(25, 9).Value = 1000
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 1000
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -826
$ws.Cells.Item(25, 14).ClearContents()
# Row 135
$ws.Cells.Item(135, 8).Value = 70000
$ws.Cells.Item(135, 10).Value = 70000
$ws.Cells.Item(135, 12).Value = 70000
$ws.Cells.Item(135, 14).Value = -80140

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 89.34614999999999
$ws.Cells.Item(2, 9).Value = 7.7647057
$ws.Cells.Item(2, 10).Value = 243.44444
$ws.Cells.Item(2, 11).Value = 46.5882342
$ws.Cells.Item(2, 12).Value = 1460.66664
$ws.Cells.Item(2, 13).Value = 66.4117658
$ws.Cells.Item(2, 14).Value = -1686.66664
# Row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()
# Row 34
$ws.Cells.Item(34, 8).Value = 960
$ws.Cells.Item(34, 9).Value = 116.666664
$ws.Cells.Item(34, 11).Value = 349.999992
$ws.Cells.Item(34, 13).Value = -265.999992
# Row 40
$ws.Cells.Item(40, 8).Value = 420.45456
$ws.Cells.Item(40, 10).Value = 495
$ws.Cells.Item(40, 12).Value = 1980
$ws.Cells.Item(40, 14).Value = -2118
# Row 46
$ws.Cells.Item(46, 8).Value = 3483.3333
$ws.Cells.Item(46, 9).Value = 300
$ws.Cells.Item(46, 10).Value = 5075
$ws.Cells.Item(46, 11).Value = 900
$ws.Cells.Item(46, 12).Value = 15225
$ws.Cells.Item(46, 13).Value = -809
$ws.Cells.Item(46, 14).Value = -15407
# Row 132
$ws.Cells.Item(132, 8).Value = 1518.3334
$ws.Cells.Item(132, 10).Value = 1518.3334
$ws.Cells.Item(132, 12).Value = 13665.0006
$ws.Cells.Item(132, 14).Value = -18725.0006
# Row 134
$ws.Cells.Item(134, 8).Value = 2000
$ws.Cells.Item(134, 9).Value = 2000
$ws.Cells.Item(134, 11).Value = 6000
$ws.Cells.Item(134, 13).Value = -930

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Cells.Item(43, 8).Value = 9286.6
$ws.Cells.Item(43, 9).Value = 811
$ws.Cells.Item(43, 10).Value = 22000
$ws.Cells.Item(43, 11).Value = 811
$ws.Cells.Item(43, 12).Value = 22000
$ws.Cells.Item(43, 13).Value = -660
$ws.Cells.Item(43, 14).Value = -22302
# Row 44
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()
# Row 102
$ws.Cells.Item(102, 8).Value = 34548.668
$ws.Cells.Item(102, 9).Value = 37544
$ws.Cells.Item(102, 11).Value = 37544
$ws.Cells.Item(102, 13).Value = -35922

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 9999
$ws.Cells.Item(16, 9).Value = 9999
$ws.Cells.Item(16, 11).Value = 9999
$ws.Cells.Item(16, 13).Value = -9829
# Row 46
$ws.Cells.Item(46, 8).Value = 255123.75
$ws.Cells.Item(46, 10).Value = 5855.7144
$ws.Cells.Item(46, 12).Value = 5855.7144
$ws.Cells.Item(46, 14).Value = -6231.7144
# Row 54
$ws.Cells.Item(54, 8).Value = 14000
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 14000
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 14000
$ws.Cells.Item(54, 13).ClearContents()
$ws.Cells.Item(54, 14).Value = -15288
# Row 93
$ws.Cells.Item(93, 8).Value = 2779
$ws.Cells.Item(93, 9).Value = 2723.75
$ws.Cells.Item(93, 10).Value = 3000
$ws.Cells.Item(93, 11).Value = 2723.75
$ws.Cells.Item(93, 12).Value = 3000
$ws.Cells.Item(93, 13).Value = -1475.75
$ws.Cells.Item(93, 14).Value = -5496
# Row 104
$ws.Cells.Item(104, 8).Value = 27670.285
$ws.Cells.Item(104, 10).Value = 27670.285
$ws.Cells.Item(104, 12).Value = 27670.285
$ws.Cells.Item(104, 14).Value = -34658.285

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Cells.Item(54, 8).Value = 27295.715
$ws.Cells.Item(54, 10).Value = 45000
$ws.Cells.Item(54, 12).Value = 45000
$ws.Cells.Item(54, 14).Value = -46040
# Row 81
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 13).ClearContents()
# Row 84
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 13).ClearContents()
